# TText - Added SIAM support, and some minor changes.
# Adds an "[Accounts]" / "Account" block (rows 20-22) to sheet1 (#TTexts),
# mirroring the existing [TTexts]/[Scopes] header+data pattern already on
# the sheet, and updates the saved selection to E22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Populate values/text (kept in original authoring order so that new
# ---- shared-string entries are appended in the same sequence as the source) ----
$ws.Range("A20").Value = "[Accounts]"
$ws.Range("B20").Value = "accUserid"
$ws.Range("C20").Value = "accPassword"
$ws.Range("D20").Value = "accPersonRef"

$ws.Range("D21").Value = "PersonRef"
$ws.Range("C21").Value = "Password"
$ws.Range("B21").Value = "UserID"
$ws.Range("A21").Value = "Account"

$ws.Range("A22").Value = "admin"
$ws.Range("B22").Value = "admin"
$ws.Range("C22").Value = "admin"

$ws.Range("D22").Value = "Ad Mini Ster"

$ws.Range("E20").Value = "autoLoginAccount"

$ws.Range("E21").Formula = "=`$A21"
$ws.Range("E22").Formula = "=`$A22"

# ---- Match formatting to the existing analogous rows (7/8 header rows, ----
# ---- 9 data row with "phrase" style, 14 data row for the E-column style) ----
$ws.Range("A7").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E7").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F20").PasteSpecial(-4122)

$ws.Range("A8").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E8").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F21").PasteSpecial(-4122)

$ws.Range("B9").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C9").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# ---- Selection recorded by Excel after editing ends on E22 ----
$ws.Range("E22").Select()
